$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting of column A style (row 10) to new rows 11 and 12
$ws.Range("A10").Copy()
$ws.Range("A11:A12").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 11
$ws.Range("A11").Value = "a1"
$ws.Range("B11").Value = 0
$ws.Range("C11").Value = [double]"1.058748410068515e-10"
$ws.Range("D11").Value = [double]"6.398161279306338e-09"
$ws.Range("E11").Value = [double]"1.648684159392424e-08"
$ws.Range("F11").Value = [double]"3.550359217045385e-24"
$ws.Range("G11").Value = [double]"1.890252129824458e-16"
$ws.Range("H11").Value = [double]"1.654609060692725e-09"
$ws.Range("I11").Value = [double]"0.2494240298614347"
$ws.Range("J11").Value = [double]"0.0001869818372195116"
$ws.Range("K11").Value = [double]"0.003904829302170598"
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = [double]"0.01662947846207714"
$ws.Range("N11").Value = [double]"99.69944323157989"
$ws.Range("O11").Value = [double]"0.01635619339710663"
$ws.Range("P11").Value = [double]"0.003318085560377239"
$ws.Range("Q11").Value = [double]"0.0007256651528104398"
$ws.Range("R11").Value = [double]"0.006362514851932103"
$ws.Range("S11").Value = [double]"0.00180134528534073"
$ws.Range("T11").Value = [double]"0.001837891659904459"
$ws.Range("U11").Value = [double]"4.137890119078464e-06"
$ws.Range("V11").Value = [double]"3.100774848934135e-06"
$ws.Range("W11").Value = [double]"5.87617926096731e-12"
$ws.Range("X11").Value = [double]"1.595264144534576e-10"
$ws.Range("Y11").Value = [double]"1.309761518072442e-14"
$ws.Range("Z11").Value = [double]"7.393762874670483e-12"
$ws.Range("AA11").Value = [double]"1.136064747761681e-12"
$ws.Range("AB11").Value = [double]"2.2827027391778e-14"
$ws.Range("AC11").Value = [double]"2.831834781864891e-13"
$ws.Range("AD11").Value = [double]"4.6586782395932e-16"
$ws.Range("AE11").Value = [double]"5.616311077334383e-16"
$ws.Range("AF11").Value = [double]"1.11233340628421e-19"
$ws.Range("AG11").Value = [double]"8.325965345192357e-21"
$ws.Range("AH11").Value = [double]"2.691280473689076e-20"
$ws.Range("AI11").Value = [double]"4.827523493279954e-21"
$ws.Range("AJ11").Value = [double]"1.087985574620105e-22"
$ws.Range("AK11").Value = 0
$ws.Range("AL11").Value = [double]"1.852026186065441e-06"
$ws.Range("AM11").Value = [double]"3.913990061825149e-07"
$ws.Range("AN11").Value = [double]"2.461398264758099e-07"
$ws.Range("AO11").Value = 0
$ws.Range("AP11").Value = 0
$ws.Range("AQ11").Value = 0
$ws.Range("AR11").Value = [double]"1.788540095378721e-16"
$ws.Range("AS11").Value = [double]"1.804963092068116e-14"
$ws.Range("AT11").Value = [double]"4.651528904454974e-14"
$ws.Range("AU11").Value = [double]"1.281085076728163e-29"
$ws.Range("AV11").Value = [double]"7.57768665315434e-22"
$ws.Range("AW11").Value = [double]"5.889935027570907e-15"
$ws.Range("AX11").Value = [double]"6.317981587697305e-10"
$ws.Range("AY11").Value = [double]"1.544600554502243e-12"
$ws.Range("AZ11").Value = [double]"3.5856195748437e-10"
$ws.Range("BA11").Value = 0
$ws.Range("BB11").Value = [double]"2.144099294670581e-07"
$ws.Range("BC11").Value = [double]"0.01992665243398281"
$ws.Range("BD11").Value = [double]"3.966321742931323e-05"
$ws.Range("BE11").Value = [double]"0.02889605237384"
$ws.Range("BF11").Value = [double]"0.03229556731351762"
$ws.Range("BG11").Value = [double]"4.041276275273107"
$ws.Range("BH11").Value = [double]"14.20581758498021"
$ws.Range("BI11").Value = [double]"29.53369384332709"
$ws.Range("BJ11").Value = [double]"12.51060501543931"
$ws.Range("BK11").Value = [double]"0.05527036064400112"
$ws.Range("BL11").Value = [double]"0.03123957934389144"
$ws.Range("BM11").Value = [double]"0.04815650032836717"
$ws.Range("BN11").Value = [double]"0.03844871584082479"
$ws.Range("BO11").Value = [double]"2.5952886020246"
$ws.Range("BP11").Value = [double]"1.360123230792329"
$ws.Range("BQ11").Value = [double]"7.322076507947115"
$ws.Range("BR11").Value = [double]"0.7641682066058152"
$ws.Range("BS11").Value = [double]"18.70049357317203"
$ws.Range("BT11").Value = [double]"0.003336991549543389"
$ws.Range("BU11").Value = [double]"2.984581412190302"
$ws.Range("BV11").Value = [double]"0.134570509131842"
$ws.Range("BW11").Value = [double]"0.008323491422675208"
$ws.Range("BX11").Value = [double]"0.05437399689901418"
$ws.Range("BY11").Value = [double]"2.006542195738793"
$ws.Range("BZ11").Value = 0
$ws.Range("CA11").Value = [double]"2.311725095296766"
$ws.Range("CB11").Value = [double]"0.8458709359497624"
$ws.Range("CC11").Value = [double]"0.3628592253619311"
$ws.Range("CD11").Value = 0
$ws.Range("CE11").Value = 0

# Row 12
$ws.Range("A12").Value = "b2"
$ws.Range("B12").Value = 0
$ws.Range("C12").Value = [double]"1.060301228771484e-10"
$ws.Range("D12").Value = [double]"6.409693707969709e-09"
$ws.Range("E12").Value = [double]"1.651596134088414e-08"
$ws.Range("F12").Value = [double]"3.495735433793975e-24"
$ws.Range("G12").Value = [double]"1.861762853801977e-16"
$ws.Range("H12").Value = [double]"1.636339142753496e-09"
$ws.Range("I12").Value = [double]"0.2494233962652553"
$ws.Range("J12").Value = [double]"0.0001869773783750931"
$ws.Range("K12").Value = [double]"0.003904822613890462"
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = [double]"0.01662943896135427"
$ws.Range("N12").Value = [double]"99.6994439962117"
$ws.Range("O12").Value = [double]"0.01635619338646314"
$ws.Range("P12").Value = [double]"0.00331816231728814"
$ws.Range("Q12").Value = [double]"0.0007256596318968826"
$ws.Range("R12").Value = [double]"0.006362428872269294"
$ws.Range("S12").Value = [double]"0.001801313301805841"
$ws.Range("T12").Value = [double]"0.001837856298155758"
$ws.Range("U12").Value = [double]"4.137814590527763e-06"
$ws.Range("V12").Value = [double]"3.100613248277515e-06"
$ws.Range("W12").Value = [double]"5.876129219956474e-12"
$ws.Range("X12").Value = [double]"1.595171001228785e-10"
$ws.Range("Y12").Value = [double]"1.309757973069216e-14"
$ws.Range("Z12").Value = [double]"7.393731422503394e-12"
$ws.Range("AA12").Value = [double]"1.13608378402341e-12"
$ws.Range("AB12").Value = [double]"2.282772605323726e-14"
$ws.Range("AC12").Value = [double]"2.83185564231592e-13"
$ws.Range("AD12").Value = [double]"4.658740640967222e-16"
$ws.Range("AE12").Value = [double]"5.615973750670056e-16"
$ws.Range("AF12").Value = [double]"1.112357195650552e-19"
$ws.Range("AG12").Value = [double]"8.326054159451756e-21"
$ws.Range("AH12").Value = [double]"2.691086836973063e-20"
$ws.Range("AI12").Value = [double]"4.827239054330904e-21"
$ws.Range("AJ12").Value = [double]"1.08800002654838e-22"
$ws.Range("AK12").Value = 0
$ws.Range("AL12").Value = [double]"1.853459167897335e-06"
$ws.Range("AM12").Value = [double]"3.916823109472563e-07"
$ws.Range("AN12").Value = [double]"2.463499727209743e-07"
$ws.Range("AO12").Value = 0
$ws.Range("AP12").Value = 0
$ws.Range("AQ12").Value = 0
$ws.Range("AR12").Value = [double]"1.792966554405978e-16"
$ws.Range("AS12").Value = [double]"1.809924530686002e-14"
$ws.Range("AT12").Value = [double]"4.663537550850302e-14"
$ws.Range("AU12").Value = [double]"1.262222281500869e-29"
$ws.Range("AV12").Value = [double]"7.46649832238714e-22"
$ws.Range("AW12").Value = [double]"5.82793348702725e-15"
$ws.Range("AX12").Value = [double]"6.318035999086779e-10"
$ws.Range("AY12").Value = [double]"1.544582134091942e-12"
$ws.Range("AZ12").Value = [double]"3.585667803374677e-10"
$ws.Range("BA12").Value = 0
$ws.Range("BB12").Value = [double]"2.144180909841319e-07"
$ws.Range("BC12").Value = [double]"0.01992698240228139"
$ws.Range("BD12").Value = [double]"3.966216004171129e-05"
$ws.Range("BE12").Value = [double]"0.02889589380848131"
$ws.Range("BF12").Value = [double]"0.03229558035182906"
$ws.Range("BG12").Value = [double]"4.041276705088187"
$ws.Range("BH12").Value = [double]"14.20581739260818"
$ws.Range("BI12").Value = [double]"29.53369356377717"
$ws.Range("BJ12").Value = [double]"12.51060495983317"
$ws.Range("BK12").Value = [double]"0.05527036444642731"
$ws.Range("BL12").Value = [double]"0.03123957926090293"
$ws.Range("BM12").Value = [double]"0.04815650021920672"
$ws.Range("BN12").Value = [double]"0.03844871575306168"
$ws.Range("BO12").Value = [double]"2.595288595901062"
$ws.Range("BP12").Value = [double]"1.360123227655681"
$ws.Range("BQ12").Value = [double]"7.32207649141538"
$ws.Range("BR12").Value = [double]"0.7641682048613712"
$ws.Range("BS12").Value = [double]"18.7004935310943"
$ws.Range("BT12").Value = [double]"0.003336991542041314"
$ws.Range("BU12").Value = [double]"2.984581405479622"
$ws.Range("BV12").Value = [double]"0.1345705088292629"
$ws.Range("BW12").Value = [double]"0.008323491403962512"
$ws.Range("BX12").Value = [double]"0.05437399677676604"
$ws.Range("BY12").Value = [double]"2.006542191227598"
$ws.Range("BZ12").Value = 0
$ws.Range("CA12").Value = [double]"2.311725090099648"
$ws.Range("CB12").Value = [double]"0.8458709340481132"
$ws.Range("CC12").Value = [double]"0.3628592245461673"
$ws.Range("CD12").Value = 0
$ws.Range("CE12").Value = 0

